$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 6284
$ws.Range("I3").Value = 6565
$ws.Range("I4").Value = 1503
$ws.Range("I5").Value = 606
$ws.Range("I6").Value = 7436
$ws.Range("I7").Value = 22394

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I7").Value = 708
$ws.Range("I8").Value = 1351
$ws.Range("I9").Value = 112
$ws.Range("I10").Value = 159
$ws.Range("I15").Value = 257
$ws.Range("I19").Value = 628
$ws.Range("I22").Value = 62
$ws.Range("I29").Value = 1368
$ws.Range("I31").Value = 225
$ws.Range("I33").Value = 1017
$ws.Range("I36").Value = 304
$ws.Range("I37").Value = 706
$ws.Range("I42").Value = 789
$ws.Range("I43").Value = 192
$ws.Range("I47").Value = 163
$ws.Range("I49").Value = 149
$ws.Range("I51").Value = 268
$ws.Range("I53").Value = 244
$ws.Range("I54").Value = 460
$ws.Range("I63").Value = 73
$ws.Range("I67").Value = 863
$ws.Range("I71").Value = 63
$ws.Range("I72").Value = 89
$ws.Range("I78").Value = 304
$ws.Range("I79").Value = 635
$ws.Range("I80").Value = 74
$ws.Range("I83").Value = 487
$ws.Range("I84").Value = 195
$ws.Range("I85").Value = 1009
$ws.Range("I86").Value = 140
$ws.Range("I88").Value = 207
$ws.Range("I89").Value = 262
$ws.Range("I90").Value = 283
$ws.Range("I94").Value = 230
$ws.Range("I96").Value = 244
$ws.Range("I97").Value = 186
$ws.Range("I99").Value = 402
$ws.Range("I101").Value = 22394

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I3").Value = 391
$ws.Range("I6").Value = 254
$ws.Range("I7").Value = 1009

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 404
$ws.Range("I3").Value = 384
$ws.Range("I6").Value = 438
$ws.Range("I7").Value = 1351

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I6").Value = 111
$ws.Range("I7").Value = 244

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 234
$ws.Range("I3").Value = 218
$ws.Range("I6").Value = 186
$ws.Range("I7").Value = 708

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I6").Value = 90
$ws.Range("I7").Value = 262

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I3").Value = 58
$ws.Range("I7").Value = 244

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 210
$ws.Range("I7").Value = 706

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I2").Value = 111
$ws.Range("I3").Value = 148
$ws.Range("I6").Value = 103
$ws.Range("I7").Value = 402

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 203
$ws.Range("I7").Value = 863

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I3").Value = 55
$ws.Range("I7").Value = 225

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I3").Value = 66
$ws.Range("I7").Value = 195

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I3").Value = 179
$ws.Range("I7").Value = 487

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I6").Value = 323
$ws.Range("I7").Value = 1017

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("I3").Value = 14
$ws.Range("I7").Value = 149

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I3").Value = 101
$ws.Range("I6").Value = 220
$ws.Range("I7").Value = 460

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 403
$ws.Range("I3").Value = 471
$ws.Range("I6").Value = 378
$ws.Range("I7").Value = 1368

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 209
$ws.Range("I6").Value = 193
$ws.Range("I7").Value = 628

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I3").Value = 244
$ws.Range("I5").Value = 26
$ws.Range("I6").Value = 269
$ws.Range("I7").Value = 789

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("I6").Value = 72
$ws.Range("I7").Value = 159

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I6").Value = 111
$ws.Range("I7").Value = 304

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 185
$ws.Range("I3").Value = 203
$ws.Range("I7").Value = 635

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I2").Value = 87
$ws.Range("I7").Value = 304

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I6").Value = 131
$ws.Range("I7").Value = 230

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I2").Value = 39
$ws.Range("I7").Value = 163

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I2").Value = 74
$ws.Range("I3").Value = 59
$ws.Range("I7").Value = 257

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("I3").Value = 40
$ws.Range("I7").Value = 112

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("I6").Value = 121
$ws.Range("I7").Value = 186

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I2").Value = 62
$ws.Range("I6").Value = 62
$ws.Range("I7").Value = 207

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I3").Value = 12
$ws.Range("I6").Value = 34
$ws.Range("I7").Value = 140

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I3").Value = 72
$ws.Range("I6").Value = 98
$ws.Range("I7").Value = 283

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I3").Value = 74
$ws.Range("I7").Value = 268

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I6").Value = 109
$ws.Range("I7").Value = 192

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("I6").Value = 17
$ws.Range("I7").Value = 62

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("I2").Value = 19
$ws.Range("I7").Value = 63

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I2").Value = 21
$ws.Range("I7").Value = 89

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("I3").Value = 15
$ws.Range("I7").Value = 74
